$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append " R" to the remarks text in P16
$ws.Range("P16").Value = "~OB Others|barcode printer issues| R"

# Update version string in C30 from 5.4.0 to 6.0.0
$ws.Range("C30").Value = "6.0.0"

# Set SICK LEAVE value for 05-13-2015 (I9) to 0.5
$ws.Range("I9").Value = 0.5

# Update TOTAL ABSENCES (I23) from 2.0 to 1.5
$ws.Range("I23").Value = 1.5
